# Add data for 2023-09-10
#
# The workbook contains year-to-date (YTD) violent-crime counts, broken
# out per calendar year (columns) and crime category (rows), both
# citywide ("Citywide Totals" / "By Neighborhood" sheets) and per
# individual Chicago neighborhood (one sheet per neighborhood). Adding a
# new day's worth of incidents nudges the YTD-through-this-calendar-day
# counts for every affected year/category/neighborhood combination up by
# the number of matching incidents recorded on 2023-09-10, and updates
# the corresponding "Total" rows/columns to match.
#
# Each block below selects a worksheet by name and rewrites the handful
# of numeric cells whose YTD counts changed (or adds a brand-new cell
# where a neighborhood previously had zero incidents of that category in
# that year).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("D2").Value = 67
$ws.Range("E2").Value = 51
$ws.Range("I2").Value = 87
$ws.Range("J3").Value = 151
$ws.Range("B6").Value = 275
$ws.Range("C6").Value = 345
$ws.Range("E6").Value = 315
$ws.Range("I6").Value = 383
$ws.Range("B7").Value = 376
$ws.Range("C7").Value = 465
$ws.Range("D7").Value = 477
$ws.Range("E7").Value = 474
$ws.Range("I7").Value = 637
$ws.Range("J7").Value = 545

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 3
$ws.Range("I7").Value = 8
$ws.Range("I28").Value = 35
$ws.Range("C32").Value = 26
$ws.Range("I32").Value = 36
$ws.Range("B36").Value = 21
$ws.Range("D47").Value = 9
$ws.Range("E50").Value = 14
$ws.Range("C53").Value = 42
$ws.Range("E53").Value = 60
$ws.Range("J65").Value = 5
$ws.Range("D70").Value = 7
$ws.Range("I70").Value = 16
$ws.Range("J74").Value = 19
$ws.Range("E78").Value = 4
$ws.Range("B98").Value = 376
$ws.Range("C98").Value = 465
$ws.Range("D98").Value = 477
$ws.Range("E98").Value = 474
$ws.Range("I98").Value = 637
$ws.Range("J98").Value = 545

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 8

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 4
$ws.Range("C6").Value = 23
$ws.Range("C7").Value = 26
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("B6").Value = 15
$ws.Range("B7").Value = 21

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("E5").Value = 12
$ws.Range("E6").Value = 14

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("E2").Value = 1
$ws.Range("E5").Value = 4

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("C6").Value = 27
$ws.Range("E6").Value = 49
$ws.Range("C7").Value = 42
$ws.Range("E7").Value = 60

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 2
$ws.Range("J6").Value = 5

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 3
$ws.Range("J6").Value = 19

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("G2").Value = 1
$ws.Range("G6").Value = 3

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("D2").Value = 1
$ws.Range("D6").Value = 9

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("D2").Value = 2
$ws.Range("I4").Value = 14
$ws.Range("D5").Value = 7
$ws.Range("I5").Value = 16
